$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns D contain price strings that Excel would otherwise auto-convert
# to floating point numbers. Temporarily mark these cells as Text so the
# literal string is preserved, then restore the General format afterward
# (this keeps the cell style identical to the original, unstyled cell).
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D10', 'D11', 'D12', 'D13', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.072.92'
$ws.Range('E2').Value = '  +1.10%  '

$ws.Range('D3').Value = '2.561.42'
$ws.Range('E3').Value = '  +2.10%  '

$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '314.91'
$ws.Range('E5').Value = '  -0.27%  '

$ws.Range('D6').Value = '98.86'
$ws.Range('E6').Value = '  +4.12%  '

$ws.Range('D7').Value = '0.574'
$ws.Range('E7').Value = '  +0.27%  '

$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('E9').Value = '  +0.76%  '

$ws.Range('D10').Value = '35.68'
$ws.Range('E10').Value = '  -0.16%  '

$ws.Range('D11').Value = '0.0815'
$ws.Range('E11').Value = '  +0.13%  '

$ws.Range('D12').Value = '7.50'
$ws.Range('E12').Value = '  -0.72%  '

$ws.Range('D13').Value = '2.984.13'
$ws.Range('E13').Value = '  +3.19%  '

$ws.Range('E14').Value = '  -0.46%  '

$ws.Range('D15').Value = '2.653.99'
$ws.Range('E15').Value = '  +5.93%  '

$ws.Range('D16').Value = '15.36'
$ws.Range('E16').Value = '  +1.68%  '

$ws.Range('D17').Value = '0.847'
$ws.Range('E17').Value = '  -0.01%  '

$ws.Range('D18').Value = '43.162.16'
$ws.Range('E18').Value = '  +1.13%  '

$ws.Range('D19').Value = '6.85'
$ws.Range('E19').Value = '  +1.43%  '

$ws.Range('D20').Value = '12.78'
$ws.Range('E20').Value = '  +0.05%  '

$ws.Range('D21').Value = '0.0₃0968'
$ws.Range('E21').Value = '  +1.08%  '

$ws.Range('D22').Value = '69.42'
$ws.Range('E22').Value = '  -0.05%  '

$ws.Range('D23').Value = '250.50'
$ws.Range('E23').Value = '  -0.03%  '

$ws.Range('E24').Value = '  +0.67%  '

$ws.Range('D25').Value = '2.12'
$ws.Range('E25').Value = '  +2.30%  '

$ws.Range('D26').Value = '27.10'
$ws.Range('E26').Value = '  +3.14%  '

$ws.Range('E27').Value = '  -0.04%  '

$ws.Range('D28').Value = '2.42'
$ws.Range('E28').Value = '  +0.15%  '

$ws.Range('D29').Value = '40.50'
$ws.Range('E29').Value = '  -2.18%  '

$ws.Range('D30').Value = '10.32'
$ws.Range('E30').Value = '  +0.64%  '

$ws.Range('D31').Value = '5.84'
$ws.Range('E31').Value = '  -1.61%  '

$ws.Range('D32').Value = '157.75'
$ws.Range('E32').Value = '  -0.76%  '

$ws.Range('D33').Value = '3.45'
$ws.Range('E33').Value = '  +6.05%  '

$ws.Range('D34').Value = '2.16'
$ws.Range('E34').Value = '  +1.83%  '

$ws.Range('E35').Value = '  +3.81%  '

$ws.Range('D36').Value = '2.68'
$ws.Range('E36').Value = '  +0.18%  '

$ws.Range('D37').Value = '18.99'
$ws.Range('E37').Value = '  -0.38%  '

$ws.Range('B38').Value = 'ApeXProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D38').Value = '2.51'
$ws.Range('E38').Value = '  +9.69%  '

$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.112'
$ws.Range('E39').Value = '  +2.08%  '

$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '24.26'
$ws.Range('E40').Value = '  +3.57%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.119'
$ws.Range('E41').Value = '  +0.63%  '

$ws.Range('D42').Value = '4.02'
$ws.Range('E42').Value = '  +7.05%  '

$ws.Range('E43').Value = '  +0.10%  '

$ws.Range('E44').Value = '  -0.29%  '

$ws.Range('D45').Value = '2.018.67'
$ws.Range('E45').Value = '  -0.21%  '

$ws.Range('D46').Value = '3.24'
$ws.Range('E46').Value = '  -2.08%  '

$ws.Range('D47').Value = '9.02'
$ws.Range('E47').Value = '  +1.36%  '

$ws.Range('D48').Value = '2.835.87'
$ws.Range('E48').Value = '  +3.17%  '

$ws.Range('D49').Value = '82.58'
$ws.Range('E49').Value = '  -2.63%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').Value = '74.91'
$ws.Range('E50').Value = '  +0.79%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.195'
$ws.Range('E51').Value = '  +3.19%  '

foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "General"
}